# welp. switched to peardeck. sorry mentimeter.
# Replace the "Example 0.2" heading's MentiMeter hyperlink reference with
# plain text pointing at Peardeck (joinpd.com).

$d = $word.ActiveDocument

$enDash = [char]0x2013
$oldHeading = "Example 0.2: Can our class speak Martian? (MentiMeter Poll)"
$newHeading = "Example 0.2: Can our class speak Martian? (Peardeck Poll " + $enDash + " joinpd.com)"

$found = $d.Content.Find.Execute($oldHeading, $false, $false, $false, $false, $false, $true, 1, $false, $newHeading, 2)

if (-not $found) {
    throw "Could not find the Example 0.2 heading text to replace."
}
